# Update "想去人数" (interested count) values in column F for sheets
# "展览" and "全部类型" to reflect newly scraped totals.

$wb = $excel.ActiveWorkbook

# Mapping of row number -> new value for column F (same updates apply to
# both the "展览" sheet and the aggregated "全部类型" sheet, since the
# latter mirrors the former's rows).
$updates = @{
    3  = 301
    4  = 1420
    5  = 8551
    9  = 256
    11 = 3447
    14 = 65
    15 = 1004
    17 = 1098
    18 = 298
    19 = 172
    20 = 2127
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
